$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply odds updates per the commit diff, grouped by row for readability.

# Row 5
$ws.Range("G5").Value = 1.91
$ws.Range("J5").Value = 2.63
$ws.Range("X5").Value = 7.5
$ws.Range("Z5").Value = 15
$ws.Range("AC5").Value = 7
$ws.Range("BD5").Value = 151

# Row 6
$ws.Range("G6").Value = 1.85
$ws.Range("I6").Value = 3.7
$ws.Range("U6").Value = 1.53
$ws.Range("V6").Value = 2.38
$ws.Range("W6").Value = 10
$ws.Range("Y6").Value = 8.5
$ws.Range("AB6").Value = 21
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 41
$ws.Range("AI6").Value = 13
$ws.Range("AM6").Value = 126
$ws.Range("AS6").Value = 101
$ws.Range("AY6").Value = 23

# Row 9
$ws.Range("G9").Value = 1.62
$ws.Range("H9").Value = 3.75
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2.2
$ws.Range("L9").Value = 5.5
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.91
$ws.Range("W9").Value = 7
$ws.Range("X9").Value = 7.5
$ws.Range("Z9").Value = 12
$ws.Range("AB9").Value = 26
$ws.Range("AD9").Value = 7
$ws.Range("AH9").Value = 29
$ws.Range("AJ9").Value = 51
$ws.Range("AK9").Value = 41
$ws.Range("AL9").Value = 41
$ws.Range("AM9").Value = 301
$ws.Range("AN9").Value = 3.5
$ws.Range("AQ9").Value = 26
$ws.Range("AU9").Value = 8.5
$ws.Range("AW9").Value = 7
$ws.Range("AX9").Value = 29
$ws.Range("AY9").Value = 34
$ws.Range("AZ9").Value = 101
$ws.Range("BA9").Value = 126
$ws.Range("BB9").Value = 251

# Row 10
$ws.Range("G10").Value = 2.45
$ws.Range("I10").Value = 3.2
$ws.Range("J10").Value = 3.2
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.5
$ws.Range("AA10").Value = 21
$ws.Range("AE10").Value = 15
$ws.Range("AH10").Value = 15
$ws.Range("AI10").Value = 12
$ws.Range("AJ10").Value = 34
$ws.Range("AK10").Value = 29
$ws.Range("AP10").Value = 26
$ws.Range("AT10").Value = 2.5
$ws.Range("AW10").Value = 5
$ws.Range("AX10").Value = 19
$ws.Range("AZ10").Value = 67
$ws.Range("BA10").Value = 101

# Row 11
$ws.Range("G11").Value = 2.15
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.88
$ws.Range("L11").Value = 4
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.4
$ws.Range("P11").Value = 3
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 1.5
$ws.Range("T11").Value = 2.5
$ws.Range("U11").Value = 1.95
$ws.Range("V11").Value = 1.8
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 21
$ws.Range("AB11").Value = 34
$ws.Range("AC11").Value = 8
$ws.Range("AD11").Value = 6
$ws.Range("AG11").Value = 9
$ws.Range("AI11").Value = 12
$ws.Range("AM11").Value = 351
$ws.Range("AT11").Value = 2.5
$ws.Range("AV11").Value = 67
$ws.Range("AW11").Value = 5

# Row 12
$ws.Range("O12").Value = 1.22
$ws.Range("P12").Value = 4.33
$ws.Range("Q12").Value = 1.73
$ws.Range("R12").Value = 2.08

# Row 19
$ws.Range("M19").Value = 1.13
$ws.Range("N19").Value = 6

# Row 20
$ws.Range("J20").Value = 16.5
$ws.Range("K20").Value = 3.55
$ws.Range("L20").Value = 1.3
$ws.Range("P20").Value = 6.9
$ws.Range("Q20").Value = 1.24
$ws.Range("R20").Value = 3.65
$ws.Range("S20").Value = 1.14
$ws.Range("T20").Value = 5
$ws.Range("U20").Value = 2.02
$ws.Range("V20").Value = 1.7
$ws.Range("AC20").Value = 28
$ws.Range("AE20").Value = 37
$ws.Range("AF20").Value = 120
$ws.Range("AG20").Value = 13
$ws.Range("AH20").Value = 8
$ws.Range("AJ20").Value = 7.1
$ws.Range("AL20").Value = 30
$ws.Range("AM20").Value = 800
$ws.Range("AN20").Value = 25
$ws.Range("AT20").Value = 5
$ws.Range("AV20").Value = 65
$ws.Range("AW20").Value = 3.6
$ws.Range("AX20").Value = 4.25
$ws.Range("AY20").Value = 12
$ws.Range("AZ20").Value = 7.2
$ws.Range("BA20").Value = 20
$ws.Range("BB20").Value = 110
